$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new numeric cells ---
$ws.Range("G1").Value = 4671
$ws.Range("I1").Value = 4672

# --- New shared strings must be created in this order so the shared-string
#     table ends up indexed as 105=Flower, 106=NONE, 107=iup ---

# "Flower" cells (rows 11-15), establishing shared string index 105
$ws.Range("H11").Value = "Flower"
$ws.Range("H12").Value = "Flower"
$ws.Range("H13").Value = "Flower"
$ws.Range("H14").Value = "Flower"
$ws.Range("H15").Value = "Flower"

# "NONE" cells (rows 3-7), establishing shared string index 106
$ws.Range("H3").Value = "NONE"
$ws.Range("H4").Value = "NONE"
$ws.Range("H5").Value = "NONE"
$ws.Range("H6").Value = "NONE"
$ws.Range("J6").Value = "NONE"
$ws.Range("H7").Value = "NONE"

# "iup" cell (row 5), establishing shared string index 107
$ws.Range("J5").Value = "iup"

# --- Remaining new numeric cells ---
$ws.Range("G3").Value = 22200
$ws.Range("G4").Value = 22300
$ws.Range("G5").Value = 22500
$ws.Range("I5").Value = 22700
$ws.Range("G6").Value = 22900
$ws.Range("I6").Value = 22900
$ws.Range("G7").Value = 23100
$ws.Range("G11").Value = 23500
$ws.Range("G12").Value = 23600
$ws.Range("G13").Value = 23800
$ws.Range("G14").Value = 24200
$ws.Range("G15").Value = 24400

# Row 7 has customFormat (style 13) applied at the row level, which would
# otherwise bleed onto newly-created cells G7/H7. Clear it so they stay
# unstyled, matching the target.
$ws.Range("G7").Style = "Normal"
$ws.Range("H7").Style = "Normal"

# --- Selection moves to J7 ---
$ws.Range("J7").Select()
